$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New residual values to be inserted into column B for rows 3-24. Values
# are parsed from strings since this PowerShell-like runtime does not
# understand scientific-notation numeric literals directly.
$newB = @{
    3  = [Convert]::ToDouble("1.052939957446597E-10")
    4  = [Convert]::ToDouble("2.480007310623478E-10")
    5  = [Convert]::ToDouble("1.503854694107076E-07")
    6  = [Convert]::ToDouble("-3.965936795080616E-07")
    7  = [Convert]::ToDouble("-3.930720193778825E-10")
    8  = [Convert]::ToDouble("1.181302580199883E-07")
    9  = [Convert]::ToDouble("2.965444589886346E-07")
    10 = [Convert]::ToDouble("3.593882045849206E-07")
    11 = [Convert]::ToDouble("-1.035781544145298E-07")
    12 = [Convert]::ToDouble("-3.87512216759589E-10")
    13 = [Convert]::ToDouble("-1.07388789361007E-07")
    14 = [Convert]::ToDouble("-1.035472805832605E-07")
    15 = [Convert]::ToDouble("6.303355340908645E-06")
    16 = [Convert]::ToDouble("-2.375649628613696E-07")
    17 = [Convert]::ToDouble("3.720025918141356E-07")
    18 = [Convert]::ToDouble("3.829984367986761E-07")
    19 = [Convert]::ToDouble("-3.160475492397508E-06")
    20 = [Convert]::ToDouble("-4.101096154340844E-08")
    21 = [Convert]::ToDouble("-1.831659499074156E-07")
    22 = [Convert]::ToDouble("2.770877186031306E-07")
    23 = [Convert]::ToDouble("2.29775004800814E-07")
    24 = [Convert]::ToDouble("-1.554241066958895E-07")
}

# Row 2: simply drop the rightmost value (K2); everything else is unchanged.
$ws.Range("K2").ClearContents()

# Rows 3-24: shift existing row values (columns B..K) one column to the
# right (dropping whatever falls off the end at column K), then write the
# new residual value into column B.
# NOTE: this runtime's `.Value` property getter does not return usable
# scalars (it returns a description string), so `.Value2` is used for
# both reads and writes instead.
for ($r = 3; $r -le 24; $r++) {
    # Capture current row values for columns B..K (col 2..11) before shifting.
    $oldValues = @{}
    for ($c = 2; $c -le 11; $c++) {
        $oldValues[$c] = $ws.Cells.Item($r, $c).Value2
    }

    # Shift each value right by one column: old column c -> new column c+1.
    # Walk from the rightmost column down so values are not clobbered
    # before being used (we already captured them all above anyway).
    for ($c = 11; $c -ge 3; $c--) {
        $srcVal = $oldValues[$c - 1]
        if ($null -eq $srcVal) {
            $ws.Cells.Item($r, $c).ClearContents()
        } else {
            $ws.Cells.Item($r, $c).Value2 = $srcVal
        }
    }

    # Insert the new residual value at column B (col 2).
    $ws.Cells.Item($r, 2).Value2 = $newB[$r]
}
